$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 8 (storemanager52 / Roma / Medhurst record).
# This shifts all subsequent rows (9-16) up by one, turning the
# former 16-row table into a 15-row table.
$ws.Rows.Item(8).Delete()

# Update the active selection to match the post-edit workbook state.
$ws.Range("H16").Select()
